# Apply updated cryptocurrency price/volume data to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings are kept as text (matches source formatting,
# e.g. "520.10" or "0.999" must not be coerced into real numbers).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = "60.340.10"
$ws.Range("E2").Value = "  -0.34%  "
$ws.Range("D3").Value = "2.620.80"
$ws.Range("E3").Value = "  +0.56%  "
$ws.Range("D5").Value = "520.10"
$ws.Range("E5").Value = "  +0.95%  "
$ws.Range("D6").Value = "150.69"
$ws.Range("E6").Value = "  -2.08%  "
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("E8").Value = "  -4.51%  "
$ws.Range("D9").Value = "6.38"
$ws.Range("E9").Value = "  -4.96%  "
$ws.Range("E10").Value = "  +1.75%  "
$ws.Range("E11").Value = "  -0.66%  "
$ws.Range("E12").Value = "  -0.98%  "
$ws.Range("D13").Value = "3.077.72"
$ws.Range("E13").Value = "  +0.50%  "
$ws.Range("D14").Value = "60.317.27"
$ws.Range("E14").Value = "  -0.39%  "
$ws.Range("D15").Value = "21.52"
$ws.Range("E15").Value = "  -0.68%  "
$ws.Range("E16").Value = "  -0.91%  "
$ws.Range("D17").Value = "2.623.66"
$ws.Range("E17").Value = "  +0.44%  "
$ws.Range("E18").Value = "  -1.99%  "
$ws.Range("D19").Value = "346.27"
$ws.Range("E19").Value = "  -3.34%  "
$ws.Range("E20").Value = "  -1.66%  "
$ws.Range("D21").Value = "6.17"
$ws.Range("E21").Value = "  -0.46%  "
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").Value = "60.91"
$ws.Range("E23").Value = "  -0.30%  "
$ws.Range("E24").Value = "  -0.88%  "
$ws.Range("E25").Value = "  -1.45%  "
$ws.Range("E26").Value = "  +0.56%  "
$ws.Range("E27").Value = "  -0.76%  "
$ws.Range("D28").Value = "7.10"
$ws.Range("E28").Value = "  -3.17%  "
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("E30").Value = "  +2.95%  "
$ws.Range("E31").Value = "  +0.45%  "
$ws.Range("D32").Value = "19.04"
$ws.Range("E32").Value = "  -2.06%  "
$ws.Range("D33").Value = "149.68"
$ws.Range("E33").Value = "  -0.52%  "
$ws.Range("E34").Value = "  -0.61%  "
$ws.Range("B35").Value = "SuiNetwork"
$ws.Range("C35").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D35").Value = "0.895"
$ws.Range("E35").Value = "  -0.19%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "1.16"
$ws.Range("E36").Value = "  -2.36%  "
$ws.Range("D38").Value = "36.65"
$ws.Range("E38").Value = "  +1.13%  "
$ws.Range("E39").Value = "  -2.07%  "
$ws.Range("E40").Value = "  -1.84%  "
$ws.Range("D41").Value = "290.20"
$ws.Range("E41").Value = "  +0.03%  "
$ws.Range("D42").Value = "0.629"
$ws.Range("E42").Value = "  +1.34%  "
$ws.Range("E43").Value = "  -1.50%  "
$ws.Range("E44").Value = "  +0.18%  "
$ws.Range("E45").Value = "  -0.86%  "
$ws.Range("D46").Value = "19.59"
$ws.Range("E46").Value = "  -0.15%  "
$ws.Range("E47").Value = "  -0.60%  "
$ws.Range("B48").Value = "WhiteBITCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D48").Value = "10.35"
$ws.Range("E48").Value = "  +0.46%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "4.74"
$ws.Range("E49").Value = "  -4.45%  "
$ws.Range("D50").Value = "18.95"
$ws.Range("E50").Value = "  -1.35%  "
$ws.Range("D51").Value = "1.963.52"
$ws.Range("E51").Value = "  -1.33%  "
